# NCI & OT update, fixed 0.9.7 bug
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("compounds")

# Open Targets Platform version: 2022.06 -> 2022.09
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2022.09"
$ws.Range("E2").Style = "Normal"

# NCI Thesaurus version: 22.08e -> 22.09d
$ws.Range("E3").Value = "22.09d"

# Update the active sheet's view/selection state
$ws.Activate()
$ws.Range("E2").Select()
